# Rename two of the "name_var" labels in the SERPAM variable dictionary:
#   "total specific cover"        -> "total specific plant cover"   (rows 20-22, V01-V03)
#   "specific vegetation cover "  -> "specific plant cover "        (rows 26-28, V07-V09)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

foreach ($r in 20..22) {
    $cell = $ws.Cells.Item($r, 2)
    if ($cell.Value() -eq "total specific cover") {
        $cell.Value = "total specific plant cover"
    }
}

foreach ($r in 26..28) {
    $cell = $ws.Cells.Item($r, 2)
    if ($cell.Value() -eq "specific vegetation cover ") {
        $cell.Value = "specific plant cover "
    }
}

# Match the saved selection state (cell B22 on Sheet1, no frozen top-left override).
$ws.Range("B22").Select()
